$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: N19 "7:00 AM START " -> "Office"
$ws.Range("N19").Value = "Office"

# Row 20: N20 "EXCEL FINANCIAL" -> "Lashaun"
$ws.Range("N20").Value = "Lashaun"

# Row 21: N21 "DARIEN CORNER MART" -> cleared
$ws.Range("N21").Value = ""

# Row 22: N22 "101 W. BELOIT ST" -> cleared
$ws.Range("N22").Value = ""

# Row 23: N23 url -> cleared
$ws.Range("N23").Value = ""

# Row 25: M25, N25, O25 cleared
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""
$ws.Range("O25").Value = ""

# Row 26: M26, N26 cleared
$ws.Range("M26").Value = ""
$ws.Range("N26").Value = ""

# Row 29: N29 "Office" -> cleared
$ws.Range("N29").Value = ""

# Row 30: N30 "Lashaun" -> cleared
$ws.Range("N30").Value = ""

# Row 37: I37, J37, K37 populated
$ws.Range("I37").Value = "#"
$ws.Range("J37").Value = "Lashaun"
$ws.Range("K37").Value = "Check In After Store"

# Row 39: J39 "Office" -> cleared
$ws.Range("J39").Value = ""

# Row 40: J40 "Kim" -> "Office"
$ws.Range("J40").Value = "Office"

# Row 41: J41 cleared -> "Kim"
$ws.Range("J41").Value = "Kim"

# New row 62
$ws.Range("A62").Value = "7)"
$ws.Range("B62").Value = "Lashaun"
$ws.Range("C62").Value = "After Store"
